# Update the "DeudoresPrueba" worksheet with the latest debtor records.
# The previous snapshot (36 rows) is replaced with an updated 30-row
# extract: some clients were removed (paid off / consolidated), two new
# clients were added ("arangos 1", "canton wok"), and several dates and
# amounts were refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing rows that no longer exist in the refreshed extract
# (the new data only spans down to row 31).
$ws.Rows("32:37").Delete()

$ws.Range("B2").Value = 'arangos 1'
$ws.Range("C2").Value = 46007
$ws.Range("D2").Value = 166000
$ws.Range("B3").Value = 'ARROZ PAISA SUBA'
$ws.Range("C3").Value = 46008
$ws.Range("D3").Value = 166000
$ws.Range("B4").Value = 'CAMILIN'
$ws.Range("C4").Value = 45997
$ws.Range("D4").Value = 166000
$ws.Range("B5").Value = 'CAMPO VERDE ZIPAQUIRA'
$ws.Range("C5").Value = 46002
$ws.Range("D5").Value = 141800
$ws.Range("B6").Value = 'canton wok'
$ws.Range("C6").Value = 46007
$ws.Range("D6").Value = 60000
$ws.Range("B7").Value = 'CIMARRON DORADO'
$ws.Range("C7").Value = 46000
$ws.Range("D7").Value = 473300
$ws.Range("B8").Value = 'COCINA CHINA'
$ws.Range("C8").Value = 46008
$ws.Range("D8").Value = 170000
$ws.Range("B9").Value = 'DARWIN FUTBOL'
$ws.Range("C9").Value = 45921
$ws.Range("D9").Value = 200000
$ws.Range("B10").Value = 'DAVIDCITO'
$ws.Range("C10").Value = 45947
$ws.Range("D10").Value = 100000
$ws.Range("B11").Value = 'EL CEBU'
$ws.Range("C11").Value = 45947
$ws.Range("D11").Value = 181800
$ws.Range("B12").Value = 'EL RUBY'
$ws.Range("C12").Value = 45992
$ws.Range("D12").Value = 85100
$ws.Range("B13").Value = 'FRANCO'
$ws.Range("C13").Value = 45996
$ws.Range("D13").Value = 20000
$ws.Range("B14").Value = 'FRESIA'
$ws.Range("C14").Value = 46006
$ws.Range("D14").Value = 248000
$ws.Range("B15").Value = 'FRIGOSOACHA'
$ws.Range("C15").Value = 46006
$ws.Range("D15").Value = 229200
$ws.Range("B16").Value = 'LA PAMPA'
$ws.Range("C16").Value = 46006
$ws.Range("D16").Value = 229900
$ws.Range("B17").Value = 'LA SELECTA'
$ws.Range("C17").Value = 45912
$ws.Range("D17").Value = 82000
$ws.Range("B18").Value = 'MERKA FRUVER DEXI'
$ws.Range("C18").Value = 45988
$ws.Range("D18").Value = 15400
$ws.Range("B19").Value = 'MERKA FRUVER DEXI'
$ws.Range("C19").Value = 45995
$ws.Range("D19").Value = 339000
$ws.Range("B20").Value = 'NEVADA'
$ws.Range("C20").Value = 46006
$ws.Range("D20").Value = 148700
$ws.Range("B21").Value = 'NOVILLON SAN MATEO'
$ws.Range("C21").Value = 45971
$ws.Range("D21").Value = 33000
$ws.Range("B22").Value = 'PARAÍSO MOSQUERA'
$ws.Range("C22").Value = 46006
$ws.Range("D22").Value = 394000
$ws.Range("B23").Value = 'PINILLA'
$ws.Range("C23").Value = 45931
$ws.Range("D23").Value = 82000
$ws.Range("B24").Value = 'PLACITA MADRILEÑA'
$ws.Range("C24").Value = 46003
$ws.Range("D24").Value = 100000
$ws.Range("B25").Value = 'PLAZA JESSICA'
$ws.Range("C25").Value = 46004
$ws.Range("D25").Value = 1705000
$ws.Range("B26").Value = 'PUNTA DE ANCA'
$ws.Range("C26").Value = 46000
$ws.Range("D26").Value = 7600
$ws.Range("B27").Value = 'SAMY 2'
$ws.Range("C27").Value = 46006
$ws.Range("D27").Value = 83000
$ws.Range("B28").Value = 'SANDRA 20 DE JULIO'
$ws.Range("C28").Value = 46000
$ws.Range("D28").Value = 300000
$ws.Range("B29").Value = 'SANTANDER SUR'
$ws.Range("C29").Value = 46006
$ws.Range("D29").Value = 56700
$ws.Range("B30").Value = 'VNZLNO PUNTA ANCA'
$ws.Range("C30").Value = 45992
$ws.Range("D30").Value = 82000
$ws.Range("B31").Value = 'WILINTONG'
$ws.Range("C31").Value = 46006
$ws.Range("D31").Value = 150000

# Keep the date format on the "Fecha" column up to date.
$ws.Range("C2:C31").NumberFormat = "yyyy\-mm\-dd"

# Re-apply the autofit column widths used for the client/date columns.
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Restore the view state (scrolled down a bit, last-used cell selected).
$ws.Range("D31").Select()
